$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "322.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.20%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.08%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "5"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.881"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "11.13%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08031"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.89%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.663"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.25%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "5"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.930"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.78%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "5"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9313"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.48%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "5"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1262"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-8.26%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "5"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1959"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.15%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "5"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8.730"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "17.62%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "5"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09203"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.34%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "5"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03531"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.63%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "5"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09554"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.05%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "5"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001297"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-7.05%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "5"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006071"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.93%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "5"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.347"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.43%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "5"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.574"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.93%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "5"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.943"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.52%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "5"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.98%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "5"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1419"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.04%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "5"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.25%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "5"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04405"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "5"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.20%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "5"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.26%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "5"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001140"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "5"
$ws.Range("B27").Value = "Spectre.aiUtilityToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("D27").Value = "--"
$ws.Range("E27").Value = "--%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "5"
$ws.Range("B28").Value = "LegolasExchange"
$ws.Range("C28").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "5"
$ws.Range("B29").Value = "BitZToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "5"
$ws.Range("B30").Value = "Birake"
$ws.Range("C30").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "5"
$ws.Range("B31").Value = "NashExchange"
$ws.Range("C31").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "5"
$ws.Range("B32").Value = "AAXToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "5"
$ws.Range("B33").Value = "CenX"
$ws.Range("C33").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "5"
$ws.Range("B34").Value = "BNIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "5"
$ws.Range("B35").Value = "UpBots"
$ws.Range("C35").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "5"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "5"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "5"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "5"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02424"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.49%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "5"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05236"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.14%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "5"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007427"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.54%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "5"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009439"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.85%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "5"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1405"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.73%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "5"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002120"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.25%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "5"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01118"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "37.52%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "5"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006731"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.07%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "5"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.05%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "5"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003001"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.87%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "5"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "5"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.05%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "5"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.05%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "5"
